# Daily update at 8 AM UTC
# Appends the new day's row of data (new "last row"), and moves the
# special "last row" date-format down from the old last row (A59) to
# the new one (A60) - mirroring the sheet's existing convention where
# the most recent day uses a plain date format while earlier days use
# the date+time format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRowFormat = $ws.Range("A59").NumberFormat
$ws.Range("A59").NumberFormat = $ws.Range("A58").NumberFormat

$ws.Range("A60").Value = 45646
$ws.Range("B60").Value = 140
$ws.Range("C60").Value = 129
$ws.Range("D60").Value = 136

$ws.Range("A60").NumberFormat = $lastRowFormat
